$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table
$cell = $tbl.Cell(2,2)
$cell.Shape.TextFrame.TextRange.Text = "changed"
